# testexcel.xlsx edit: "up to 95, testing ALL even more"
#
# Sheet-name -> physical-file mapping (via workbook relationships):
#   input-0 -> sheet1.xml
#   input-1 -> sheet2.xml
#   input-2 -> sheet3.xml

$wb = $excel.ActiveWorkbook

# NOTE: these contain U+00A0 (NBSP) between the words, matching the
# original workbook's shared-string text exactly (not regular spaces).
$successMsg = "a TORPEDO parancs eredménye SUCCESS"
$failMsg    = "a TORPEDO parancs eredménye FAIL"

# ---------------------------------------------------------------
# input-1 (sheet2.xml): a new TORPEDO,ALL step block appended
#   - existing last step (row 19/20) flips from SINGLE/FAIL to ALL/SUCCESS
#   - a fresh SINGLE/FAIL block is appended as rows 22-23
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("input-1")

$ws1.Range("B19").Value = "TORPEDO,ALL"
$ws1.Range("B20").Value = $successMsg

$ws1.Range("A22").Value = "A teszt lépése"
$ws1.Range("B22").Value = "TORPEDO,SINGLE"
$ws1.Range("A23").Value = "Elvárt kimenet/eredmény"
$ws1.Range("B23").Value = $failMsg

# ---------------------------------------------------------------
# input-2 (sheet3.xml): preconditions tweak (1 -> 2 secondary torpedoes)
#   plus a new TORPEDO,ALL / FAIL block appended as rows 13-14
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("input-2")

$ws2.Range("B2").Value = "primary-ben 0, secondaryben 2 torpedo"

$ws2.Range("A13").Value = "A teszt lépése"
$ws2.Range("B13").Value = "TORPEDO,ALL"
$ws2.Range("A14").Value = "Elvárt kimenet/eredmény"
$ws2.Range("B14").Value = $failMsg

# ---------------------------------------------------------------
# View state: input-2 was the active/selected tab before, input-1 is now.
# Apply input-2's new selection first so activating input-1 last leaves it
# as the saved ActiveTab, matching the recorded final state.
# ---------------------------------------------------------------
[void]$ws2.Range("D6").Select()

[void]$ws1.Activate()
[void]$ws1.Range("D21").Select()
